# The workbook "Metadata" sheet (sheet1) lists ValueSet properties as
# Property/Value pairs, one per row. This change:
#   1. Inserts a new "Jurisdiction" property row (with an empty value)
#      right after the "Contact" row, pushing Description/Purpose/
#      Copyright/Immutable down by one row.
#   2. Refreshes the "Date" property value to the new publication date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 11 (right after "Contact"/row 10),
# pushing the existing rows 11-14 down to 12-15.
$ws.Rows.Item(11).Insert()

# Populate the new row with the "Jurisdiction" property and an empty value.
# Setting the value to a single apostrophe forces Excel to store it as an
# explicit (empty) text value rather than leaving the cell completely blank,
# matching how the other text properties in this sheet are stored.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "'"

# Re-apply the formatting used by the rest of the data rows (border/
# alignment) to the freshly inserted row, copying it from the row right
# below (which held the formatting previously used at row 11).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the "Date" property (row 8, column B) to the new timestamp.
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
